# Auto-generated edit script: updates Leve profit tables (H-N columns)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ==== ALC sheet ====
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 3856.9312
$ws.Range("I100").Value = 1570.75
$ws.Range("K100").Value = 1570.75
$ws.Range("M100").Value = -1029.75

# Row 113
$ws.Range("H113").Value = 2219
$ws.Range("I113").Value = 2684.6667
$ws.Range("J113").Value = 1869.75
$ws.Range("K113").Value = 2684.6667
$ws.Range("L113").Value = 1869.75
$ws.Range("M113").Value = 569.3332999999998
$ws.Range("N113").Value = -8377.75

# Row 124
$ws.Range("H124").Value = 63750
$ws.Range("J124").Value = 63750
$ws.Range("L124").Value = 63750
$ws.Range("N124").Value = -73570

# Row 137
$ws.Range("H137").Value = 2462.2856
$ws.Range("I137").Value = 1648.875
$ws.Range("K137").Value = 4946.625
$ws.Range("M137").Value = -2396.625

# Row 138
$ws.Range("H138").Value = 3415.611
$ws.Range("I138").Value = 3169.25
$ws.Range("K138").Value = 9507.75
$ws.Range("M138").Value = -4367.75

# ==== ARM sheet ====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3477.29
$ws.Range("I32").Value = 1920.0115
$ws.Range("J32").Value = 13899.077
$ws.Range("K32").Value = 1920.0115
$ws.Range("L32").Value = 13899.077
$ws.Range("M32").Value = -1633.0115
$ws.Range("N32").Value = -14473.077

# Row 64
$ws.Range("H64").Value = 59999.668

# Row 67
$ws.Range("H67").Value = 59999.668

# Row 74
$ws.Range("H74").Value = 2085799.9
$ws.Range("J74").Value = 5513
$ws.Range("L74").Value = 5513
$ws.Range("N74").Value = -7261

# Row 77
$ws.Range("H77").Value = 2085799.9
$ws.Range("J77").Value = 5513
$ws.Range("L77").Value = 27565
$ws.Range("N77").Value = -36301

# Row 102
$ws.Range("H102").Value = 4199.4443
$ws.Range("I102").Value = 4199.4443
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4199.4443
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2577.4443
$ws.Range("N102").ClearContents()

# ==== BSM sheet ====
$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 100181
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 100181
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 105
$ws.Range("H105").Value = 3954.0605
$ws.Range("J105").Value = 4546.8
$ws.Range("L105").Value = 4546.8
$ws.Range("N105").Value = -8040.8

# ==== CRP sheet ====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9910.200000000001
$ws.Range("I31").Value = 3068.5
$ws.Range("J31").Value = 11620.625
$ws.Range("K31").Value = 3068.5
$ws.Range("L31").Value = 11620.625
$ws.Range("M31").Value = -2773.5
$ws.Range("N31").Value = -12210.625

# Row 34
$ws.Range("H34").Value = 9910.200000000001
$ws.Range("I34").Value = 3068.5
$ws.Range("J34").Value = 11620.625
$ws.Range("K34").Value = 3068.5
$ws.Range("L34").Value = 11620.625
$ws.Range("M34").Value = -2866.5
$ws.Range("N34").Value = -12024.625

# Row 58
$ws.Range("H58").Value = 774328.2
$ws.Range("I58").Value = 1124634.1
$ws.Range("J58").Value = 3655
$ws.Range("K58").Value = 1124634.1
$ws.Range("L58").Value = 3655
$ws.Range("M58").Value = -1124431.1
$ws.Range("N58").Value = -4061

# Row 88
$ws.Range("H88").Value = 15268.2
$ws.Range("J88").Value = 15268.2
$ws.Range("L88").Value = 15268.2
$ws.Range("N88").Value = -16080.2

# Row 91
$ws.Range("H91").Value = 15268.2
$ws.Range("J91").Value = 15268.2
$ws.Range("L91").Value = 15268.2
$ws.Range("N91").Value = -18076.2

# Row 114
$ws.Range("H114").Value = 100684
$ws.Range("J114").Value = 100684
$ws.Range("L114").Value = 100684
$ws.Range("N114").Value = -109362

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

# Row 119
$ws.Range("H119").Value = 100761
$ws.Range("J119").Value = 100761
$ws.Range("L119").Value = 100761
$ws.Range("N119").Value = -110437

# Row 120
$ws.Range("H120").Value = 54634
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 54634
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 54634
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -61892

# Row 136
$ws.Range("H136").Value = 774328.2
$ws.Range("I136").Value = 1124634.1
$ws.Range("J136").Value = 3655
$ws.Range("K136").Value = 3373902.3
$ws.Range("L136").Value = 10965
$ws.Range("M136").Value = -3371352.3
$ws.Range("N136").Value = -16065

# ==== CUL sheet ====
$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 4997.8
$ws.Range("I70").Value = 2999.6667
$ws.Range("K70").Value = 8999.000100000001
$ws.Range("M70").Value = -8684.000100000001

# Row 73
$ws.Range("H73").Value = 4997.8
$ws.Range("I73").Value = 2999.6667
$ws.Range("K73").Value = 8999.000100000001
$ws.Range("M73").Value = -7907.000100000001

# Row 98
$ws.Range("H98").Value = 813.1667
$ws.Range("J98").Value = 778
$ws.Range("L98").Value = 2334
$ws.Range("N98").Value = -5330

# ==== GSM sheet ====
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 208.59091
$ws.Range("I2").Value = 64.38461
$ws.Range("J2").Value = 416.8889
$ws.Range("K2").Value = 64.38461
$ws.Range("L2").Value = 416.8889
$ws.Range("M2").Value = 48.61539
$ws.Range("N2").Value = -642.8888999999999

# Row 33
$ws.Range("H33").Value = 186666
$ws.Range("I33").Value = 500000
$ws.Range("K33").Value = 500000
$ws.Range("M33").Value = -499748

# Row 132
$ws.Range("H132").Value = 2183.923
$ws.Range("I132").Value = 2079.3
$ws.Range("J132").Value = 2532.6667
$ws.Range("K132").Value = 6237.900000000001
$ws.Range("L132").Value = 7598.000100000001
$ws.Range("M132").Value = -3707.900000000001
$ws.Range("N132").Value = -12658.0001

# Row 135
$ws.Range("H135").Value = 74576.664
$ws.Range("J135").Value = 74576.664
$ws.Range("L135").Value = 74576.664
$ws.Range("N135").Value = -84716.664

# ==== LTW sheet ====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2377.6
$ws.Range("I40").Value = 2377.6
$ws.Range("K40").Value = 2377.6
$ws.Range("M40").Value = -2241.6

# Row 136
$ws.Range("H136").Value = 4205.7905
$ws.Range("I136").Value = 3568.8918
$ws.Range("K136").Value = 10706.6754
$ws.Range("M136").Value = -8156.6754

# ==== WVR sheet ====
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 996.881
$ws.Range("I107").Value = 724.4138
$ws.Range("J107").Value = 1604.6923
$ws.Range("K107").Value = 2173.2414
$ws.Range("L107").Value = 4814.0769
$ws.Range("M107").Value = -253.2413999999999
$ws.Range("N107").Value = -8654.0769

# Row 132
$ws.Range("H132").Value = 3608.279
$ws.Range("I132").Value = 3739.4
$ws.Range("K132").Value = 11218.2
$ws.Range("M132").Value = -8688.200000000001

